# Weekly data update: a new week's price record for Piña (Vega Monumental
# Concepción) is inserted at row 289, pushing the existing rows 289-313
# down to 290-314 (dimension grows from A1:T313 to A1:T314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 289 (existing rows 289-313 shift to 290-314)
$ws.Rows.Item(289).Insert()

# The newly inserted row is blank; populate it with the same record as the
# row that used to occupy position 289 (now at row 290 after the shift),
# then overwrite the handful of fields that differ for the new week.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(289, $col).Value = $ws.Cells.Item(290, $col).Value2
}

# Fields that change for the new week's entry
$ws.Cells.Item(289, 4).Value  = 45223   # D289 Fecha
$ws.Cells.Item(289, 14).Value = 20000   # N289 Precio minimo
$ws.Cells.Item(289, 15).Value = 21000   # O289 Precio maximo
$ws.Cells.Item(289, 16).Value = 20500   # P289 Precio promedio ponderado
$ws.Cells.Item(289, 19).Value = 1464    # S289 Precio $/Kg
